$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 584 (pushes existing 584.. down to 586..)
$ws.Range("A584:A585").EntireRow.Insert()

# Fill in row 584 (new)
$ws.Cells.Item(584, 1).Value = 6
$ws.Cells.Item(584, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(584, 3).Value = "Metropolitana"
$ws.Cells.Item(584, 4).Value = 44939
$ws.Cells.Item(584, 5).Value = 13
$ws.Cells.Item(584, 6).Value = 100112052
$ws.Cells.Item(584, 7).Value = "Albahaca"
$ws.Cells.Item(584, 8).Value = "Sin especificar"
$ws.Cells.Item(584, 9).Value = "Primera"
$ws.Cells.Item(584, 10).Value = 990
$ws.Cells.Item(584, 11).Value = 2500
$ws.Cells.Item(584, 12).Value = 3000
$ws.Cells.Item(584, 13).Value = 2652
$ws.Cells.Item(584, 14).Value = "`$/docena de matas"
$ws.Cells.Item(584, 15).Value = "Región Metropolitana"
$ws.Cells.Item(584, 16).Value = 442
$ws.Cells.Item(584, 17).Value = 6
$ws.Cells.Item(584, 18).Value = "Hortaliza"

# Fill in row 585 (new)
$ws.Cells.Item(585, 1).Value = 6
$ws.Cells.Item(585, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(585, 3).Value = "Metropolitana"
$ws.Cells.Item(585, 4).Value = 44939
$ws.Cells.Item(585, 5).Value = 13
$ws.Cells.Item(585, 6).Value = 100112052
$ws.Cells.Item(585, 7).Value = "Albahaca"
$ws.Cells.Item(585, 8).Value = "Sin especificar"
$ws.Cells.Item(585, 9).Value = "Segunda"
$ws.Cells.Item(585, 10).Value = 200
$ws.Cells.Item(585, 11).Value = 2000
$ws.Cells.Item(585, 12).Value = 2000
$ws.Cells.Item(585, 13).Value = 2000
$ws.Cells.Item(585, 14).Value = "`$/docena de matas"
$ws.Cells.Item(585, 15).Value = "Región Metropolitana"
$ws.Cells.Item(585, 16).Value = 333
$ws.Cells.Item(585, 17).Value = 6
$ws.Cells.Item(585, 18).Value = "Hortaliza"
